$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-04-09 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-10 Wednesday", 2)
$d.Content.Find.Execute("58×32=", $true, $false, $false, $false, $false, $true, 1, $false, "17×93=", 2)
$d.Content.Find.Execute("14×21=", $true, $false, $false, $false, $false, $true, 1, $false, "47×82=", 2)
$d.Content.Find.Execute("97×77=", $true, $false, $false, $false, $false, $true, 1, $false, "89×91=", 2)
$d.Content.Find.Execute("25×98=", $true, $false, $false, $false, $false, $true, 1, $false, "65×67=", 2)
$d.Content.Find.Execute("43×17=", $true, $false, $false, $false, $false, $true, 1, $false, "23×18=", 2)
$d.Content.Find.Execute("20×83=", $true, $false, $false, $false, $false, $true, 1, $false, "33×19=", 2)
$d.Content.Find.Execute("75×71=", $true, $false, $false, $false, $false, $true, 1, $false, "25×25=", 2)
$d.Content.Find.Execute("82×27=", $true, $false, $false, $false, $false, $true, 1, $false, "75×45=", 2)
$d.Content.Find.Execute("76×66=", $true, $false, $false, $false, $false, $true, 1, $false, "33×97=", 2)
$d.Content.Find.Execute("61×97=", $true, $false, $false, $false, $false, $true, 1, $false, "27×55=", 2)
$d.Content.Find.Execute("49×61=", $true, $false, $false, $false, $false, $true, 1, $false, "35×38=", 2)
$d.Content.Find.Execute("41×64=", $true, $false, $false, $false, $false, $true, 1, $false, "13×67=", 2)
$d.Content.Find.Execute("54×22=", $true, $false, $false, $false, $false, $true, 1, $false, "13×81=", 2)
$d.Content.Find.Execute("91×27=", $true, $false, $false, $false, $false, $true, 1, $false, "51×21=", 2)
$d.Content.Find.Execute("80×90=", $true, $false, $false, $false, $false, $true, 1, $false, "41×65=", 2)
$d.Content.Find.Execute("23×76=", $true, $false, $false, $false, $false, $true, 1, $false, "60×82=", 2)
$d.Content.Find.Execute("33×60=", $true, $false, $false, $false, $false, $true, 1, $false, "39×48=", 2)
$d.Content.Find.Execute("67×14=", $true, $false, $false, $false, $false, $true, 1, $false, "35×74=", 2)
$d.Content.Find.Execute("73×96=", $true, $false, $false, $false, $false, $true, 1, $false, "65×80=", 2)
$d.Content.Find.Execute("31×89=", $true, $false, $false, $false, $false, $true, 1, $false, "49×36=", 2)
$d.Content.Find.Execute("53×55=", $true, $false, $false, $false, $false, $true, 1, $false, "54×94=", 2)
$d.Content.Find.Execute("52×84=", $true, $false, $false, $false, $false, $true, 1, $false, "82×14=", 2)
$d.Content.Find.Execute("36×58=", $true, $false, $false, $false, $false, $true, 1, $false, "89×16=", 2)
$d.Content.Find.Execute("56×32=", $true, $false, $false, $false, $false, $true, 1, $false, "19×55=", 2)
$d.Content.Find.Execute("69×49=", $true, $false, $false, $false, $false, $true, 1, $false, "68×97=", 2)
